$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the shared-string table in the same order as the target workbook:
# "...completed" must land before "...continued" in xl/sharedStrings.xml,
# even though "continued" is used on the earlier row (56).
$ws.Range("D57").Value = "Manual: 2nd review completed"
$ws.Range("D56").Value = "Manual: 2nd review continued"

# Row 56: 2012-12-11, 1h effort
$ws.Range("A56").Value = 41254
$ws.Range("B56").Value = 1

# Row 57: 2012-12-12, 1.75h effort
$ws.Range("A57").Value = 41255
$ws.Range("B57").Value = 1.75

# Copy the date format (style index 1, ddd dd/mm/yyyy) from the row above
# down onto the two new date cells without introducing new style records.
$ws.Range("A55").Copy()
$ws.Range("A56:A57").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D56").Select()
